# Update res_bus/vm_pu.xlsx results for the "case with 380 kV" re-run.
# Bus 1 (column B) setpoint moved from 1.05 pu to 1.02 pu, and the resulting
# voltage magnitudes for buses 2-5 and 7-11 (columns C,D,E,F,I,J,K,L,M) were
# recomputed for all 24 time steps (rows 2-25). Columns A, G, H and N are
# unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B,C,D,E,F,I,J,K,L,M (skipping G=1 (unchanged) and H=empty)
$colIndices = @(2,3,4,5,6,9,10,11,12,13)

$data = @(
    @(1.02,1.080416947809461,1.082600923681428,1.083516436617712,1.093851892338386,1.06169137221202,1.085296588463696,1.085268716931678,1.086181849270658,1.096490727932794),
    @(1.02,1.081638820709447,1.083590410843401,1.084587261377796,1.094955990059692,1.062087799923795,1.086178520799054,1.086076453293009,1.087070899070883,1.097414893440482),
    @(1.02,1.08242949720248,1.084230661042608,1.085280408757779,1.095670700494448,1.062343137207978,1.086748650882024,1.086598477465925,1.087645820267207,1.098012574447814),
    @(1.02,1.082761909246758,1.084499819287961,1.08557186882606,1.095971233636998,1.062450199306736,1.086988205182483,1.086817784879367,1.087887433080355,1.098263764481661),
    @(1.02,1.082817723413898,1.084545011982563,1.085620809841621,1.096021698529141,1.062468158994483,1.087028419894741,1.08685459868778,1.087927996022098,1.098305936008285),
    @(1.02,1.082433938861022,1.084234257559232,1.085284303021457,1.095674715961672,1.06234456888346,1.086751852321928,1.086601408457126,1.087649049039021,1.098015931155577),
    @(1.02,1.080829877378409,1.082935328704264,1.083878274975154,1.094224968576908,1.061825590901965,1.085594753659824,1.085541827074231,1.086482381341756,1.096803119733041),
    @(1.02,1.07800358515331,1.080646321009911,1.081402580407307,1.091672488852452,1.060902049681848,1.083551638393821,1.083669819673398,1.084423831009781,1.094663546456152),
    @(1.02,1.076119484360376,1.079120197503145,1.079753361865139,1.089972250070598,1.060280259538648,1.082186719757439,1.082418491065269,1.08304958857978,1.093235482703581),
    @(1.02,1.075303650135023,1.078459333213261,1.079039518038702,1.089236353995436,1.060009566565048,1.081595010112788,1.081875854864657,1.082454072334562,1.092616707253954),
    @(1.02,1.075000610225723,1.078213851588009,1.078774405635405,1.088963056371648,1.059908800229038,1.081375117985725,1.081674174179076,1.082232801052845,1.092386803344553),
    @(1.02,1.075065613468905,1.078266508568452,1.078831271305859,1.089021677530015,1.059930424864547,1.081422290353154,1.081717440879809,1.082280267677366,1.092436121385875),
    @(1.02,1.07527860082832,1.078439041777002,1.079017602933969,1.089213762176227,1.060001241656954,1.081576835905169,1.08185918635852,1.0824357834152,1.092597704618628),
    @(1.02,1.075409828950438,1.078545344203241,1.079132413403036,1.089332118059629,1.060044845196821,1.081672042600233,1.08194650430878,1.082531592487219,1.092697252985656),
    @(1.02,1.076173628305319,1.079164055922363,1.079800743124409,1.090021095702849,1.060298193868044,1.082225974934593,1.082454487059196,1.083089101259544,1.093276539984574),
    @(1.02,1.076652736523997,1.079552145443572,1.080220042684818,1.090453357813935,1.060456723232165,1.082573256365617,1.082772915724421,1.083438687943874,1.093639799831212),
    @(1.02,1.076932191771469,1.079778507424008,1.080464639973406,1.090705519785528,1.060549050497295,1.082775753158739,1.082958572389826,1.083642551599131,1.093851643228288),
    @(1.02,1.077027478764914,1.079855690342808,1.080548045925972,1.090791505752574,1.060580507938558,1.082844788068531,1.083021863414496,1.083712056343034,1.093923869607292),
    @(1.02,1.076601332808011,1.079510507537382,1.080175053028997,1.090406977001387,1.060439729039994,1.08253600324965,1.082738759352766,1.083401185194002,1.093600829623505),
    @(1.02,1.075215881450222,1.07838823525248,1.07896273178966,1.089157196746938,1.059980393923773,1.08153132899379,1.081817449207109,1.082389989846951,1.092550124134547),
    @(1.02,1.074344776360513,1.077682576062069,1.078200733952209,1.088371681450417,1.059690324892975,1.080899043455632,1.081237481918057,1.081753805567564,1.091889138044787),
    @(1.02,1.07480656785666,1.078056663598814,1.07860466129039,1.088788072574325,1.059844216222074,1.081234287857264,1.081545000474658,1.082091097638277,1.09223957425406),
    @(1.02,1.076624559935571,1.079529321927203,1.080195381834329,1.09042793439867,1.060447408417046,1.08255283653464,1.082754193385703,1.083418131209096,1.093618438705976),
    @(1.02,1.078734226506662,1.081238102912341,1.082042384705401,1.092332113952471,1.061141880158449,1.084080330242829,1.084154361628211,1.084956343393022,1.095216971038559)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $colIndices.Length; $j++) {
        $ws.Cells.Item($excelRow, $colIndices[$j]).Value = $rowValues[$j]
    }
}

Write-Output "Updated $($data.Length) rows"